$wb = $excel.ActiveWorkbook

# --- 1. Set selection on Lead_Config5 (Rel_Based_Asgmt / Existing Account row)
#        without making it the active sheet (matches target: selection added,
#        no tabSelected) -------------------------------------------------
$s5 = $wb.Worksheets.Item("Lead_Config5")
$s5.Range("A2:G2").Select()

# --- 2. Set selection on Lead_Config15 (the sheet that currently owns
#        tabSelected) to the full data range; tabSelected will be dropped
#        automatically once a different sheet becomes active later --------
$s15 = $wb.Worksheets.Item("Lead_Config15")
$s15.Range("A1:G2").Select()

# --- 3. Add the new Lead_Config16 worksheet at the end of the tab strip ---
$s16 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$s16.Name = "Lead_Config16"

# Copy layout/format (header styling, border styling) from Lead_Config15,
# the most similar existing sheet, then fix up the values that differ for
# this new config row.
$s15.Range("A1:G2").Copy($s16.Range("A1"))

# Match the source sheet's column widths too.
for ($i = 1; $i -le 7; $i++) {
    $s16.Columns.Item($i).ColumnWidth = $s15.Columns.Item($i).ColumnWidth
}

$s16.Range("E1").Value = "Lead Fields"
$s16.Range("C2").Value = "Lead Fields"
$s16.Range("F2").Value = "N.A."
$s16.Range("G2").Value = "Checked"

# Leave the new sheet active, selecting G2 to match the authored selection.
$s16.Range("G2").Select()
